$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, pushing old rows 164..273 down to 165..274.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new record.
# Columns A-C, E-I, K-M, P-R keep the same values the (now shifted) row had,
# i.e. the same as the data that is now in row 165 (originally row 164).
$ws.Range("A164").Value = 5
$ws.Range("B164").Value = 'Macroferia Regional de Talca'
$ws.Range("C164").Value = 'Maule'
$ws.Range("D164").Value = 45216
$ws.Range("E164").Value = 7
$ws.Range("F164").Value = 100112031
$ws.Range("G164").Value = 'Poroto verde'
$ws.Range("H164").Value = 'Sin especificar'
$ws.Range("I164").Value = 'Primera'
$ws.Range("J164").Value = 150
$ws.Range("K164").Value = 30000
$ws.Range("L164").Value = 30000
$ws.Range("M164").Value = 30000
$ws.Range("N164").Value = '$/malla 25 kilos'
$ws.Range("O164").Value = 'Perú'
$ws.Range("P164").Value = 1200
$ws.Range("Q164").Value = 25
$ws.Range("R164").Value = 'Hortaliza'
